$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1 (16:30 -> 16:45)
$ws.Range("F1").Value = "Last status check on: 25.02.2022 16:45"

# Row 7 (MOL Olomoucka): convert delta (D7) and old-date (E7) from text to real numbers
$ws.Range("D7").Value = 1

$ws.Range("E7").Value = 44617.68914351852
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
